# Insert 9 new data rows at row 766 (pushing existing rows 766-816 down to 775-825)
# and populate them with the new Femacal de La Calera - Nectarin price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank rows before the current row 766, shifting everything below down.
$ws.Rows.Item(766).Resize(9).Insert()

# Data for the newly inserted rows (row, K-Variedad, L-Calidad, M-Volumen, N-PrecioMin,
# O-PrecioMax, P-PrecioProm, Q-Unidad, R-Origen, S-Precio$/Kg, T-Kg/unidad).
# Columns A,B,C,D,E,F,G,H,I,J share the same values across all rows in this sheet.
$newRows = @(
    @{ Row=766; K="August Red";   L="Especial";                M=85; N=16000; O=16000; P=16000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=1067; T=15 },
    @{ Row=767; K="August Red";   L="Primera";                 M=87; N=14000; O=14000; P=14000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=933;  T=15 },
    @{ Row=768; K="August Red";   L="Segunda";                 M=80; N=12000; O=12000; P=12000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=800;  T=15 },
    @{ Row=769; K="August pearl"; L="Especial";                M=65; N=16000; O=16000; P=16000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=1067; T=15 },
    @{ Row=770; K="August pearl"; L="Extra (doble especial)";  M=60; N=18000; O=18000; P=18000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=1200; T=15 },
    @{ Row=771; K="August pearl"; L="Primera";                 M=60; N=14000; O=14000; P=14000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=933;  T=15 },
    @{ Row=772; K="Venus";        L="Especial";                M=70; N=16000; O=16000; P=16000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=1067; T=15 },
    @{ Row=773; K="Venus";        L="Primera";                 M=78; N=14000; O=14000; P=14000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=933;  T=15 },
    @{ Row=774; K="Venus";        L="Segunda";                 M=75; N=12000; O=12000; P=12000; Q="`$/caja 15 kilos empedrada"; R="Región de O'Higgins"; S=800;  T=15 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 44610
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100103
    $ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($row, 9).Value = 100103006
    $ws.Cells.Item($row, 10).Value = "Nectarín"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
